# Update "Forecast Comparison" sheet with corrected forecast output:
#  - Insert a new "Week_Start_Date" column after "Week"
#  - Shorten week labels ("W01" -> "W1", etc.)
#  - Populate the new Week_Start_Date column with the week's start date (as text)
#  - Correct a handful of MyForecast values
#  - Store is_holiday_week as a boolean instead of a number
# Also refresh the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# Helper: assign a value to a cell while keeping it plain text, even if the
# text looks like a number or a date (Excel would otherwise silently
# reinterpret "2025-01-05" as a date serial, or "2394" as a number).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new, empty column B ("Week_Start_Date"); everything from the old
# column B onward shifts one column to the right (B->C, C->D, ... I->J).
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Week_Start_Date"

# New, corrected data per week: Week label, Week start date, MyForecast.
# (Amazon Mean/P70/P80/P90 forecasts, Product Title and is_holiday_week are
# unaffected in value, only shifted right by the column insert above.)
$weekData = @(
    @{ Row = 2;  Week = "W1";  Date = "2025-01-05"; MyForecast = 183 },
    @{ Row = 3;  Week = "W2";  Date = "2025-01-12"; MyForecast = 169 },
    @{ Row = 4;  Week = "W3";  Date = "2025-01-19"; MyForecast = 161 },
    @{ Row = 5;  Week = "W4";  Date = "2025-01-26"; MyForecast = 162 },
    @{ Row = 6;  Week = "W5";  Date = "2025-02-02"; MyForecast = 174 },
    @{ Row = 7;  Week = "W6";  Date = "2025-02-09"; MyForecast = 182 },
    @{ Row = 8;  Week = "W7";  Date = "2025-02-16"; MyForecast = 140 },
    @{ Row = 9;  Week = "W8";  Date = "2025-02-23"; MyForecast = 138 },
    @{ Row = 10; Week = "W9";  Date = "2025-03-02"; MyForecast = 166 },
    @{ Row = 11; Week = "W10"; Date = "2025-03-09"; MyForecast = 163 },
    @{ Row = 12; Week = "W11"; Date = "2025-03-16"; MyForecast = 157 },
    @{ Row = 13; Week = "W12"; Date = "2025-03-23"; MyForecast = 149 },
    @{ Row = 14; Week = "W13"; Date = "2025-03-30"; MyForecast = 130 },
    @{ Row = 15; Week = "W14"; Date = "2025-04-06"; MyForecast = 112 },
    @{ Row = 16; Week = "W15"; Date = "2025-04-13"; MyForecast = 103 },
    @{ Row = 17; Week = "W16"; Date = "2025-04-20"; MyForecast = 105 }
)

foreach ($entry in $weekData) {
    $row = $entry.Row

    # Week label (column A)
    $ws.Cells.Item($row, 1).Value = $entry.Week

    # Week_Start_Date (column B) - keep it a plain text value, not an
    # auto-converted date serial number.
    Set-TextValue $ws.Cells.Item($row, 2) $entry.Date

    # MyForecast (column D, after the insert)
    $ws.Cells.Item($row, 4).Value = $entry.MyForecast

    # is_holiday_week (column J, after the insert) becomes a boolean.
    $ws.Cells.Item($row, 10).Value = $false
}

# Update the Summary sheet totals to reflect the corrected forecast values.
# These are stored as text in the sheet, so keep them text here too.
$summary = $wb.Worksheets.Item("Summary")
Set-TextValue $summary.Range("B9")  "2394"
Set-TextValue $summary.Range("B10") "1309"
Set-TextValue $summary.Range("B11") "675"
